$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 for columns P and Q (copy format from O1, which
# already carries the bold/bordered header style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing data rows 2-25: swap I<->K style updates and M<->O style updates,
# and add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new = 2
}
